$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record is being inserted before the current last data row (row 18).
# Insert a new row at 18, which pushes the existing row 18 down to row 19 (unchanged).
$ws.Rows.Item(18).Insert()

# The new row 18 gets the data that row 17 held previously (row 17 is about to be
# updated with a newer observation), i.e. it is effectively a copy of the old row 17.
$ws.Cells.Item(18, 1).Value2 = 8
$ws.Cells.Item(18, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(18, 3).Value2 = "Coquimbo"
$ws.Cells.Item(18, 4).Value2 = 44827
$ws.Cells.Item(18, 5).Value2 = 4
$ws.Cells.Item(18, 6).Value2 = 100112039
$ws.Cells.Item(18, 7).Value2 = "Ciboulette"
$ws.Cells.Item(18, 8).Value2 = "Sin especificar"
$ws.Cells.Item(18, 9).Value2 = "Primera"
$ws.Cells.Item(18, 10).Value2 = 1200
$ws.Cells.Item(18, 11).Value2 = 2000
$ws.Cells.Item(18, 12).Value2 = 2500
$ws.Cells.Item(18, 13).Value2 = 2250
$ws.Cells.Item(18, 14).Value2 = "`$/docena de atados"
$ws.Cells.Item(18, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(18, 16).Value2 = 750
$ws.Cells.Item(18, 17).Value2 = 3
$ws.Cells.Item(18, 18).Value2 = "Hortaliza"

# Row 17 is updated with the newest observation (new date and updated volume).
$ws.Cells.Item(17, 4).Value2 = 45006
$ws.Cells.Item(17, 10).Value2 = 1100
